# "tablas y evaluaciones todas"
# Re-zeroes the res_c1 (A), res_c3 (C) and total (E) columns for every
# entity row (2-71), since those scores were computed from data that no
# longer applies; only row 72 keeps (new, recomputed) non-zero figures,
# along with fresh per-question (p*) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-71: reset columns A (res_c1), C (res_c3) and E (total) to 0.
for ($r = 2; $r -le 71; $r++) {
    $ws.Range("A$r").Value = 0
    $ws.Range("C$r").Value = 0
    $ws.Range("E$r").Value = 0
}

# Row 23 also has its per-question answers (I, M, O, Q..AF) cleared to 0.
$row23Cols = @("I","M","O","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF")
foreach ($col in $row23Cols) {
    $ws.Range("$col" + "23").Value = 0
}

# Row 72 (totals/last entity row) gets the newly recomputed figures.
$ws.Range("A72").Value = 22.77
$ws.Range("B72").Value = 17.64
$ws.Range("C72").Value = 17.97
$ws.Range("D72").Value = 6.3
$ws.Range("E72").Value = 64.68000000000001
$ws.Range("U72").Value = 1.57
$ws.Range("V72").Value = 1.09
$ws.Range("Y72").Value = 1.1
$ws.Range("AC72").Value = 3.1
$ws.Range("AD72").Value = 1.1
$ws.Range("AF72").Value = 1.1
$ws.Range("AG72").Value = 13.47
$ws.Range("AH72").Value = 2.2
$ws.Range("AI72").Value = 2.3
$ws.Range("AK72").Value = 2.1
$ws.Range("AR72").Value = 4.2
